$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 109 — this shifts the existing rows 109:134 down to 110:135,
# preserving all of their data/formatting, matching the dimension change
# (A1:R134 -> A1:R135).
$ws.Rows(109).Insert()

# Populate the newly inserted row 109 with the new record.
$ws.Cells.Item(109, 1).Value = 5
$ws.Cells.Item(109, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(109, 3).Value = "Maule"
$ws.Cells.Item(109, 4).Value = 44543
$ws.Cells.Item(109, 5).Value = 7
$ws.Cells.Item(109, 6).Value = 100112024
$ws.Cells.Item(109, 7).Value = "Choclo"
$ws.Cells.Item(109, 8).Value = "Choclero"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 20000
$ws.Cells.Item(109, 11).Value = 350
$ws.Cells.Item(109, 12).Value = 350
$ws.Cells.Item(109, 13).Value = 350
$ws.Cells.Item(109, 14).Value = "$/unidad"
$ws.Cells.Item(109, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(109, 16).Value = 350
$ws.Cells.Item(109, 17).Value = 1
$ws.Cells.Item(109, 18).Value = "Hortaliza"
